$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record arrived (newest first); insert it as row 4,
# pushing the existing rows 4-12 down to 5-13.
$ws.Rows.Item(4).Insert()

$ws.Cells.Item(4, 1).Value = 7
$ws.Cells.Item(4, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(4, 3).Value = "Ñuble"
$ws.Cells.Item(4, 4).Value = 45152
$ws.Cells.Item(4, 5).Value = 16
$ws.Cells.Item(4, 6).Value = "Fruta"
$ws.Cells.Item(4, 7).Value = 100102
$ws.Cells.Item(4, 8).Value = "Cítricos"
$ws.Cells.Item(4, 9).Value = 100102006
$ws.Cells.Item(4, 10).Value = "Pomelo"
$ws.Cells.Item(4, 11).Value = "Start Ruby"
$ws.Cells.Item(4, 12).Value = "Primera"
$ws.Cells.Item(4, 13).Value = 60
$ws.Cells.Item(4, 14).Value = 16000
$ws.Cells.Item(4, 15).Value = 16000
$ws.Cells.Item(4, 16).Value = 16000
$ws.Cells.Item(4, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(4, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(4, 19).Value = 1143
$ws.Cells.Item(4, 20).Value = 14
